$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2000, 'Seven Star Residency', 21.245918, 72.861221),
    @(7438, 'VNSGU', 21.153851, 72.779714999999996),
    @(473, 'Sarthana Nature Park', 21.229827, 72.899540999999999),
    @(8493, 'Maharana Partap Garden', 21.224121, 72.884878),
    @(9843, 'Adajan Patiya', 21.198701, 72.807055000000005),
    @(8372, 'Gujrat Gas Circle', 21.193932, 72.801903999999993),
    @(6582, 'Athva Gate', 21.185621000000001, 72.810660999999996),
    @(7461, 'Majura Gate', 21.180719, 72.818521000000004),
    @(4781, 'Surat Railway Station', 21.204922, 72.839980999999995),
    @(3213, 'Gotalavadi', 21.211634, 72.830855999999997),
    @(8664, 'Kiran Hospital', 21.219275, 72.837095000000005),
    @(9823, 'Vallabhacharya RD', 21.219588000000002, 72.851616000000007),
    @(9503, 'Hirabaug Circle', 21.215689000000001, 72.862960999999999),
    @(8872, 'Kapodra Fire Station', 21.220329, 72.874052000000006),
    @(2032, 'Yogi Nagar Brts', 21.217689, 72.895975000000007),
    @(7463, 'Simada Canal Junc.', 21.210167999999999, 72.896876000000006)
)

$row = 7
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).Value = $item[3]
    $row++
}

$ws.Range("A9").Select() | Out-Null
